# Add two new columns "I0" (I) and "IF" (J) to the sheet.
# I0 is always 1; IF is a copy of the existing "IP" column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, same formatting (bold / border / centered) as the
# other header cells, obtained by copying the adjacent "IP" header (H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill the data rows: I = 1, J = same value as column H ("IP") for that row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, "H").End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
